$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert two new rows for DistilBERT "Oversampling" entries.
#    New rows become rows 10 and 11 (inheriting format from row 9
#    above them, which already carries the right styles/number formats).
# ------------------------------------------------------------------
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(10).Insert()

# ------------------------------------------------------------------
# 2) Row 3, 4, 5: add "Others" (L) column note.
# ------------------------------------------------------------------
$ws.Range("L3").Value = "1 hidden layer"
$ws.Range("L4").Value = "1 hidden layer"
$ws.Range("L5").Value = "1 hidden layer"

# ------------------------------------------------------------------
# 3) Row 9: new NonToxicScore DistilBERT Upsampling 512 tokens entry.
#    Mark C9 with the same green highlight used on C8, add the model
#    file name, and extend the "Others" note to mention 512 tokens.
# ------------------------------------------------------------------
$ws.Range("C9").Interior.Color = 5296274
$ws.Range("K9").Value = "DistilBertToxicClassification512tok.pth"
$ws.Range("L9").Value = "2 hidden layers, , GELU(), 2 epochs. 512 tokens"

# ------------------------------------------------------------------
# 4) Fill the two newly inserted rows (10 and 11) with the new
#    DistilBERT / Oversampling results.
# ------------------------------------------------------------------
$ws.Range("A10").Value = "DistilBERT"
$ws.Range("B10").Value = 1
$ws.Range("D10").Value = 64
$ws.Range("E10").Value = 0.00002
$ws.Range("F10").Value = 16
$ws.Range("G10").Value = "256/ 32"
$ws.Range("H10").Value = 0.1
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = "Oversampling"
$ws.Range("K10").Value = "DistilBertToxicClassification6.pth"
$ws.Range("L10").Value = "2 hidden layers, , GELU(), 2 epochs"

$ws.Range("A11").Value = "DistilBERT"
$ws.Range("B11").Value = 1
$ws.Range("D11").Value = 64
$ws.Range("E11").Value = 0.00002
$ws.Range("F11").Value = 16
$ws.Range("G11").Value = "256/ 32"
$ws.Range("H11").Value = 0.1
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = "Oversampling"
$ws.Range("K11").Value = "DistilBertToxicClassification7.pth"
$ws.Range("L11").Value = "2 hidden layers, , GELU(), 2 epochs. 512 tokens"

# ------------------------------------------------------------------
# 5) Rows 12-14 (previously 10-12, the BERT rows) now need an
#    "Others" note too; row 13's note also changes wording.
# ------------------------------------------------------------------
$ws.Range("L12").Value = "1 hidden layer"
$ws.Range("L13").Value = "1 hidden layer. scheduler Tmax = epochs = 4"
$ws.Range("L14").Value = "1 hidden layer"

# ------------------------------------------------------------------
# 6) Hidden Size (G) column: right align header + data (new cellXfs).
# ------------------------------------------------------------------
$ws.Range("G1:G1048576").HorizontalAlignment = -4152

# ------------------------------------------------------------------
# 7) Column K width grew to fit the longer file names.
# ------------------------------------------------------------------
$ws.Columns.Item(11).ColumnWidth = 35.85546875

# ------------------------------------------------------------------
# 8) Dimension / AutoFilter / defined name all need to grow to L14.
# ------------------------------------------------------------------
$ws.AutoFilter.Range.AutoFilter(1)
$ws.Range("A1:L14").AutoFilter(1)
$names = $wb.Names
$n = $names.Item(1)
$n.RefersTo = "=Sheet1!`$A`$1:`$L`$14"

# ------------------------------------------------------------------
# 9) Selection used by Excel when the file was saved.
# ------------------------------------------------------------------
$ws.Range("K19").Select()
